$wb = $excel.ActiveWorkbook

# --- Step 1: the existing "effort" sheet becomes "Effort R 0.9" --------------
$oldSheet = $wb.Worksheets.Item(1)
$oldSheet.Name = "Effort R 0.9"

# --- Step 2: insert a brand new sheet in front of it for the new revision ----
$newSheet = $wb.Worksheets.Add($oldSheet)
$newSheet.Name = "Effort R 1.0"

# Re-resolve stable references by name (Add() re-indexes the collection).
$rev09 = $wb.Worksheets.Item("Effort R 0.9")
$rev10 = $wb.Worksheets.Item("Effort R 1.0")

# --- Step 3: populate "Effort R 1.0" with the new effort log -----------------
$rev10.Range("A1").Value = "Date"
$rev10.Range("B1").Value = "Effort [h]"
$rev10.Range("C1").Value = "Additional Effort [h]"
$rev10.Range("D1").Value = "Task"

$data = @(
    @(41423, 2.5,  2.5,  "Revision of Makefile"),
    @(41424, 2,    $null, "Documentation Makefile changes. Concept for new sync objects"),
    @(41425, 2,    $null, "Concept of new sync objects"),
    @(41426, 0.75, $null, "Design of implementation new sync objects "),
    @(41430, 2,    $null, "Implementation of mutexes"),
    @(41431, 2,    $null, "Implementation of mutexes"),
    @(41432, 2,    2,     "Update Manual"),
    @(41432, 2.25, $null, "Implementation of mutexes"),
    @(41435, 2,    $null, "Implementation of mutexes: Basically done. No test case implemented yet, no testing done yet"),
    @(41436, 1.5,  $null, "Implementation of tc11_mutex"),
    @(41439, 1.5,  2.5,  "Implementation of semaphores and first, very preliminary but successfuls tests"),
    @(41440, 2.5,  $null, "Implementation tc12_queue"),
    @(41442, 2,    $null, "Implementation tc12_queue"),
    @(41443, 1,    3,    "Implementation tc12_queue"),
    @(41444, 1.5,  2.5,  "Design and implementation tc13_eventStates")
)

$r = 2
foreach ($row in $data) {
    $rev10.Cells.Item($r, 1).Value = $row[0]
    $rev10.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $rev10.Cells.Item($r, 3).Value = $row[2]
    }
    $rev10.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Date formatting for the whole "Date" column, header included
# (matches style s="1" used throughout column A on sheet "Effort R 0.9").
$rev10.Range("A1:A16").NumberFormat = "ddd\ dd/mm/yyyy"

# --- Step 4: cosmetics matching the authored sheet ---------------------------
$rev10.Columns.Item(1).ColumnWidth = 13.45
$rev10.Columns.Item(4).ColumnWidth = 71.74

$rev10.PageSetup.PaperSize = 9
$rev10.PageSetup.Orientation = 1

# --- Step 5: selections (order matters: last Select() wins the active tab) --
$rev09.Range("B1:B1048576").Select()
$rev10.Range("D25").Select()
